# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Reorders the "Periodo Mora" column (E16:E22) so the most recent period
# (2308) is listed first and the oldest (2302) last, carrying the
# "Valor Mora" (column F) figures along with their period as the table
# is re-sorted (the 40000 value stays tied to period 2308, 46400 to the
# rest).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order for column E (Periodo Mora), rows 16..22, newest period first
$periods = @("2308", "2307", "2306", "2305", "2304", "2303", "2302")

# Corresponding Valor Mora (column F) for each period, in the same order
$valores = @(40000, 46400, 46400, 46400, 46400, 46400, 46400)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
